$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right count 6 -> 9, Wrong count 3 -> 2
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 ("Total"): Right total 168 -> 252, and matching "Total/Max" label text
$ws.Range("B12").Value = 252
$ws.Range("E12").Value = "252/252"
